$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values per the data refresh
$ws.Range("A3").Value = 10003
$ws.Range("A25").Value = 10003

# Update view: scroll to top-left A1 (remove topLeftCell override) and select full rows 34 through end
$ws.Range("A34:XFD1048576").Select()
